# Insert a new column "elev" before the existing "distance" column (column G)
# and populate it, per the commit:
#   "Added elevation field to station, fixed flag to be a string not an int."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at G; everything from G onward shifts right by one.
$ws.Columns("G:G").Insert()

# New header cell for the inserted column.
$ws.Range("G2").Value = "elev"

# Populate elevation value for each station row (rows 4-21).
$elev = 5.0999999999999996
for ($r = 4; $r -le 21; $r++) {
    $ws.Cells.Item($r, 7).Value = $elev
}

# Match the cursor position left behind in the saved workbook.
[void]$ws.Range("G24").Select()
